# Remove the trailing "footer" block that was scraped from the Jupiter
# course-catalog page: a blank paragraph, the "Ver no Jupiter ..." line
# and the "(c) 2020 ... Creative Commons Attribution" line. The
# paragraph right after the (now removed) block, and the one right
# before it (the "1996.OMETTO..." bibliography line), are left intact.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph
# dynamically (rather than hard-coding an index) by scanning the
# document's Paragraphs collection.
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text.StartsWith("Ver no Jupiter Salvar em pdf Salvar em docx")) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 2) {
    $prevPara = $d.Paragraphs.Item($targetIndex - 1)
    $nextPara = $d.Paragraphs.Item($targetIndex + 1)

    $prevText = $prevPara.Range.Text.Trim()
    $nextText = $nextPara.Range.Text

    # Safety check: the paragraph before should be blank, and the one
    # after should be the copyright/footer line (identified by a
    # distinctive substring so we don't depend on the exact encoding
    # of the leading "(c)" glyph).
    if ($prevText.Length -eq 0 -and $nextText.Contains("luizeleno@usp.br")) {
        $startPos = $prevPara.Range.Start
        $endPos = $nextPara.Range.End
        $deleteRange = $d.Range($startPos, $endPos)
        $deleteRange.Delete()
    }
}
